{"js": "// Replace the date and all the two-digit division answers in the table.\n// Every old value is unique in the document, so a plain search & replace\n// (matching whole, case-sensitive strings) for each pair is safe.\nconst replacements = [\n  [\"2023-10-15 Sunday\", \"2023-10-16 Monday\"],\n  [\"24\u00f73=8, 0\", \"39\u00f78=4, 7\"],\n  [\"35\u00f75=7, 0\", \"92\u00f75=18, 2\"],\n  [\"63\u00f74=15, 3\", \"42\u00f75=8, 2\"],\n  [\"84\u00f74=21, 0\", \"52\u00f78=6, 4\"],\n  [\"33\u00f75=6, 3\", \"22\u00f73=7, 1\"],\n  [\"84\u00f73=28, 0\", \"28\u00f75=5, 3\"],\n  [\"18\u00f76=3, 0\", \"52\u00f73=17, 1\"],\n  [\"28\u00f79=3, 1\", \"62\u00f73=20, 2\"],\n  [\"58\u00f77=8, 2\", \"11\u00f73=3, 2\"],\n  [\"61\u00f74=15, 1\", \"37\u00f72=18, 1\"],\n  [\"76\u00f78=9, 4\", \"93\u00f73=31, 0\"],\n  [\"44\u00f76=7, 2\", \"83\u00f72=41, 1\"],\n  [\"89\u00f76=14, 5\", \"67\u00f77=9, 4\"],\n  [\"40\u00f74=10, 0\", \"72\u00f72=36, 0\"],\n  [\"15\u00f73=5, 0\", \"57\u00f75=11, 2\"],\n  [\"23\u00f78=2, 7\", \"86\u00f73=28, 2\"],\n  [\"73\u00f79=8, 1\", \"88\u00f79=9, 7\"],\n  [\"82\u00f78=10, 2\", \"18\u00f78=2, 2\"],\n  [\"38\u00f77=5, 3\", \"94\u00f76=15, 4\"],\n  [\"42\u00f73=14, 0\", \"61\u00f77=8, 5\"],\n  [\"23\u00f75=4, 3\", \"22\u00f78=2, 6\"],\n  [\"38\u00f72=19, 0\", \"63\u00f77=9, 0\"],\n  [\"93\u00f72=46, 1\", \"35\u00f76=5, 5\"],\n  [\"11\u00f79=1, 2\", \"62\u00f74=15, 2\"],\n  [\"57\u00f76=9, 3\", \"91\u00f74=22, 3\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date and all the two-digit division answers in the table.\n# Every old value is unique in the document, so a plain Find/Replace\n# (whole-document range, case-sensitive, replace-all) for each pair is safe.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2023-10-15 Sunday\", \"2023-10-16 Monday\"),\n    @(\"24\u00f73=8, 0\", \"39\u00f78=4, 7\"),\n    @(\"35\u00f75=7, 0\", \"92\u00f75=18, 2\"),\n    @(\"63\u00f74=15, 3\", \"42\u00f75=8, 2\"),\n    @(\"84\u00f74=21, 0\", \"52\u00f78=6, 4\"),\n    @(\"33\u00f75=6, 3\", \"22\u00f73=7, 1\"),\n    @(\"84\u00f73=28, 0\", \"28\u00f75=5, 3\"),\n    @(\"18\u00f76=3, 0\", \"52\u00f73=17, 1\"),\n    @(\"28\u00f79=3, 1\", \"62\u00f73=20, 2\"),\n    @(\"58\u00f77=8, 2\", \"11\u00f73=3, 2\"),\n    @(\"61\u00f74=15, 1\", \"37\u00f72=18, 1\"),\n    @(\"76\u00f78=9, 4\", \"93\u00f73=31, 0\"),\n    @(\"44\u00f76=7, 2\", \"83\u00f72=41, 1\"),\n    @(\"89\u00f76=14, 5\", \"67\u00f77=9, 4\"),\n    @(\"40\u00f74=10, 0\", \"72\u00f72=36, 0\"),\n    @(\"15\u00f73=5, 0\", \"57\u00f75=11, 2\"),\n    @(\"23\u00f78=2, 7\", \"86\u00f73=28, 2\"),\n    @(\"73\u00f79=8, 1\", \"88\u00f79=9, 7\"),\n    @(\"82\u00f78=10, 2\", \"18\u00f78=2, 2\"),\n    @(\"38\u00f77=5, 3\", \"94\u00f76=15, 4\"),\n    @(\"42\u00f73=14, 0\", \"61\u00f77=8, 5\"),\n    @(\"23\u00f75=4, 3\", \"22\u00f78=2, 6\"),\n    @(\"38\u00f72=19, 0\", \"63\u00f77=9, 0\"),\n    @(\"93\u00f72=46, 1\", \"35\u00f76=5, 5\"),\n    @(\"11\u00f79=1, 2\", \"62\u00f74=15, 2\"),\n    @(\"57\u00f76=9, 3\", \"91\u00f74=22, 3\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
